# Add team W/L/T record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (AD1:AF1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, top-aligned, thin border)
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1        # xlContinuous

# Fill team record for every data row (2-58) with the season's W/L/T totals
$ws.Range("AD2:AD58").Value = 69
$ws.Range("AE2:AE58").Value = 93
$ws.Range("AF2:AF58").Value = 0
